$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.026244282722473
$ws.Range("B1").Value = 1.694387316703796
$ws.Range("C1").Value = 4.407351016998291
$ws.Range("D1").Value = 2.441533327102661
$ws.Range("E1").Value = 1.339480519294739
